$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 2-10 (data rows) following the updated NATMI computation
# (Ligand-expressing cells / Receptor-expressing cells count changed 1 -> 3,
#  with derived statistics recalculated accordingly).

# NOTE: PowerShell hashtables are case-insensitive, so a key named "r" would
# collide with a key named "R" (the column-R value). Use "RowNum" for the
# row index to avoid any collision with column letters.
$rows = @(
    @{ RowNum = 2;  E = 3; G = 7.636417666666667;  H = 22.909253;          I = 0.108532481296676;  J = 0.108532481296676;  K = 3; M = 9.682562333333333; N = 29.047687;          O = 0.5358521175370563; P = 0.5358521175370563; Q = 73.94009006086789;  R = 665.460810547811;   S = 0.05815735992437478; T = 0.05815735992437479 },
    @{ RowNum = 3;  E = 3; G = 7.636417666666667;  H = 22.909253;          I = 0.108532481296676;  J = 0.108532481296676;  K = 3; M = 4.196433666666667; N = 12.589301;          O = 0.2322389248810544; P = 0.2322389248810544; Q = 32.04572018912811;  R = 288.411481702153;   S = 0.02520546677101317; T = 0.02520546677101318 },
    @{ RowNum = 4;  E = 3; G = 7.636417666666667;  H = 22.909253;          I = 0.108532481296676;  J = 0.108532481296676;  K = 3; M = 4.190471333333334; N = 12.571414;          O = 0.2319089575818893; P = 0.2319089575818892; Q = 32.00018932152689;  R = 288.001703893742;   S = 0.02516965460128802; T = 0.02516965460128802 },
    @{ RowNum = 5;  E = 3; G = 15.103385;          H = 45.31015499999999;  I = 0.214656652056136;  J = 0.214656652056136;  K = 3; M = 9.682562333333333; N = 29.047687;          O = 0.5358521175370563; P = 0.5358521175370563; Q = 146.2394667068316;  R = 1316.155200361485;  S = 0.1150242215476956;  T = 0.1150242215476956 },
    @{ RowNum = 6;  E = 3; G = 15.103385;          H = 45.31015499999999;  I = 0.214656652056136;  J = 0.214656652056136;  K = 3; M = 4.196433666666667; N = 12.589301;          O = 0.2322389248810544; P = 0.2322389248810544; Q = 63.38035329462832;  R = 570.423179651655;   S = 0.04985163009208358; T = 0.0498516300920836 },
    @{ RowNum = 7;  E = 3; G = 15.103385;          H = 45.31015499999999;  I = 0.214656652056136;  J = 0.214656652056136;  K = 3; M = 4.190471333333334; N = 12.571414;          O = 0.2319089575818893; P = 0.2319089575818892; Q = 63.29030187879666;  R = 569.61271690917;    S = 0.0497808004163568;  T = 0.0497808004163568 },
    @{ RowNum = 8;  E = 3; G = 47.62086333333334;  H = 142.86259;          I = 0.676810866647188;  J = 0.676810866647188;  K = 3; M = 9.682562333333333; N = 29.047687;          O = 0.5358521175370563; P = 0.5358521175370563; Q = 461.0919775921478;  R = 4149.82779832933;   S = 0.3626705360649859;  T = 0.3626705360649859 },
    @{ RowNum = 9;  E = 3; G = 47.62086333333334;  H = 142.86259;          I = 0.676810866647188;  J = 0.676810866647188;  K = 3; M = 4.196433666666667; N = 12.589301;          O = 0.2322389248810544; P = 0.2322389248810544; Q = 199.8377941277323;  R = 1798.54014714959;   S = 0.1571818280179576;  T = 0.1571818280179576 },
    @{ RowNum = 10; E = 3; G = 47.62086333333334;  H = 142.86259;          I = 0.676810866647188;  J = 0.676810866647188;  K = 3; M = 4.190471333333334; N = 12.571414;          O = 0.2319089575818893; P = 0.2319089575818892; Q = 199.5538626669178;  R = 1795.98476400226;   S = 0.1569585025642444;  T = 0.1569585025642444 }
)

foreach ($row in $rows) {
    $rn = $row.RowNum
    $ws.Range("E$rn").Value = $row.E
    $ws.Range("G$rn").Value = $row.G
    $ws.Range("H$rn").Value = $row.H
    $ws.Range("I$rn").Value = $row.I
    $ws.Range("J$rn").Value = $row.J
    $ws.Range("K$rn").Value = $row.K
    $ws.Range("M$rn").Value = $row.M
    $ws.Range("N$rn").Value = $row.N
    $ws.Range("O$rn").Value = $row.O
    $ws.Range("P$rn").Value = $row.P
    $ws.Range("Q$rn").Value = $row.Q
    $ws.Range("R$rn").Value = $row.R
    $ws.Range("S$rn").Value = $row.S
    $ws.Range("T$rn").Value = $row.T
}

$wb.Save()
